$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Each entry: row, Coin (B), Link (C), Price (D), Volume1h (E)
$data = @(
    ,@(2, 'Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '29.115.08', '  -2.39%  ')
    ,@(3, 'Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '1.850.05', '  -1.46%  ')
    ,@(4, 'TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '1.000', '  -0.19%  ')
    ,@(5, 'XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.6949', '  -4.68%  ')
    ,@(6, 'BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '238.91', '  -1.25%  ')
    ,@(7, 'USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '1.000', '  -0.22%  ')
    ,@(8, 'Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.3044', '  -2.89%  ')
    ,@(9, 'Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.07581', '  +6.81%  ')
    ,@(10, 'Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '23.38', '  -4.15%  ')
    ,@(11, 'TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.08131', '  -1.75%  ')
    ,@(12, 'Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '0.7261', '  -2.90%  ')
    ,@(13, 'Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '5.230', '  -1.89%  ')
    ,@(14, 'WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '1.808.17', '  -4.95%  ')
    ,@(15, 'Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '89.37', '  -3.47%  ')
    ,@(16, 'WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '29.105.80', '  -2.61%  ')
    ,@(17, 'Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '5.781', '  -4.57%  ')
    ,@(18, 'ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.000007781', '  -0.67%  ')
    ,@(19, 'Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '13.13', '  -1.88%  ')
    ,@(20, 'BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '237.07', '  -4.49%  ')
    ,@(21, 'Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '0.9994', '  -0.41%  ')
    ,@(22, 'WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '2.103.42', '  -3.43%  ')
    ,@(23, 'BinanceUSD', 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd', '1.000', '  -0.16%  ')
    ,@(24, 'Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '7.593', '  -2.01%  ')
    ,@(25, 'Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '9.014', '  -1.78%  ')
    ,@(26, 'Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '160.77', '  -1.40%  ')
    ,@(27, 'Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.1454', '  -5.71%  ')
    ,@(28, 'EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '18.10', '  -2.49%  ')
    ,@(29, 'LidoDAOToken', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', '1.984', '  -2.17%  ')
    ,@(30, 'Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '1.394', '  -3.16%  ')
    ,@(31, 'Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '4.490', '  -1.38%  ')
    ,@(32, 'PancakeSwap', 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake', '1.492', '  -2.42%  ')
    ,@(33, 'InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '4.001', '  -4.58%  ')
    ,@(34, 'Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.05214', '  -0.92%  ')
    ,@(35, 'ARBITRUM', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', '1.190', '  -3.62%  ')
    ,@(36, 'Frax', 'https://coinranking.com/coin/KfWtaeV1W+frax-frax', '1.035', '  +3.54%  ')
    ,@(37, 'ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '0.7035', '  -7.16%  ')
    ,@(38, 'HuobiToken', 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht', '2.654', '  -2.25%  ')
    ,@(39, 'VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.01860', '  -3.72%  ')
    ,@(40, 'MXToken', 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx', '2.692', '  -2.26%  ')
    ,@(41, 'TrustWalletToken', 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt', '0.9373', '  +7.83%  ')
    ,@(42, 'FraxShare', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', '6.032', '  +0.48%  ')
    ,@(43, 'Maker', 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr', '1.072.69', '  +0.31%  ')
    ,@(44, 'TheSandbox', 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand', '0.4273', '  -4.90%  ')
    ,@(45, 'Aave', 'https://coinranking.com/coin/ixgUfzmLR+aave-aave', '70.14', '  -1.74%  ')
    ,@(46, 'PaxDollar', 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp', '0.9999', '  -0.21%  ')
    ,@(47, 'Quant', 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt', '103.25', '  -1.53%  ')
    ,@(48, 'RenderToken', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', '1.777', '  -3.15%  ')
    ,@(49, 'RocketPoolETH', 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth', '1.995.11', '  -3.71%  ')
    ,@(50, 'EnergySwap', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', '9.173', '  -3.59%  ')
    ,@(51, 'Aptos', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', '7.043', '  -6.31%  ')
)

foreach ($item in $data) {
    $row = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]

    # Force Price and Volume columns to be stored as text so values
    # like '1.000' or '5.230' are not re-interpreted as numbers.
    $ws.Cells.Item($row, 4).NumberFormat = "@"
    $ws.Cells.Item($row, 4).Value = $item[3]
    $ws.Cells.Item($row, 5).NumberFormat = "@"
    $ws.Cells.Item($row, 5).Value = $item[4]
}

Write-Output "Updated $($data.Count) rows"